# Apply gradebook updates: new homework/exam scores entered for several
# students, plus a selection/cursor move.
#
# xlPasteFormats constant used with PasteSpecial to copy only the cell
# *formatting* (fill/border/font/number format) from a template cell that
# already carries the desired style, without disturbing the destination
# cell's value/formula.
$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-Format($srcAddress, $dstAddress) {
    $src = $ws.Range($srcAddress)
    $dst = $ws.Range($dstAddress)
    $src.Copy()
    $dst.PasteSpecial($xlPasteFormats)
}

# --- Row 10: Волнухин Михаил -------------------------------------------
# Homework cells D10:L10 lose their "missing" (green) highlight and go back
# to the plain/default look; the student's variant (M10) is recorded as 3,
# which drives the variant-tally formulas in N10:Q10 (formatting there is
# refreshed too).
Copy-Format "N27" "D10:L10"
Copy-Format "N27" "N10"
$ws.Range("M10").Value = 3

# --- Row 14: exam score ---------------------------------------------------
Copy-Format "S16" "S14"
$ws.Range("S14").Value = 5

# --- Row 17: homework scores filled in ------------------------------------
Copy-Format "J34" "E17:L17"
$ws.Range("E17:L17").Value = 5

# --- Row 18: remaining homework + exam score ------------------------------
Copy-Format "J34" "K18:L18"
$ws.Range("K18:L18").Value = 5
$ws.Range("S18").Value = 5

# --- Row 21: homework scores + exam score ---------------------------------
Copy-Format "E3" "E21"
Copy-Format "E3" "K21:L21"
$ws.Range("K21:L21").Value = 5
$ws.Range("S21").Value = 5

# --- Row 33: homework scores -----------------------------------------------
Copy-Format "E3" "F33"
$ws.Range("F33").Value = 5
Copy-Format "E3" "H33:L33"
$ws.Range("H33:L33").Value = 5

# --- Move the active selection to match the saved cursor position --------
$ws.Range("S14").Select()
